$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.217.71'
$ws.Range('E2').Value = '  -4.63%  '
$ws.Range('D3').Value = '2.573.48'
$ws.Range('E3').Value = '  -3.71%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'507.29"
$ws.Range('D6').Value = "'144.83"
$ws.Range('E6').Value = '  -7.12%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = "'0.572"
$ws.Range('D9').Value = '2.589.27'
$ws.Range('E9').Value = '  -3.74%  '
$ws.Range('D10').Value = "'6.29"
$ws.Range('E10').Value = '  -4.20%  '
$ws.Range('E11').Value = '  -4.97%  '
$ws.Range('E12').Value = '  -5.31%  '
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('D14').Value = '3.025.10'
$ws.Range('E14').Value = '  -3.71%  '
$ws.Range('D15').Value = '58.220.23'
$ws.Range('E15').Value = '  -4.61%  '
$ws.Range('D16').Value = "'21.06"
$ws.Range('E16').Value = '  -4.70%  '
$ws.Range('E17').Value = '  -4.28%  '
$ws.Range('D18').Value = '2.578.50'
$ws.Range('E18').Value = '  -3.93%  '
$ws.Range('E19').Value = '  -5.34%  '
$ws.Range('D20').Value = "'342.46"
$ws.Range('E20').Value = '  -3.64%  '
$ws.Range('D21').Value = "'10.31"
$ws.Range('E21').Value = '  -3.94%  '
$ws.Range('D22').Value = "'6.06"
$ws.Range('E22').Value = '  -4.53%  '
$ws.Range('D23').Value = "'0.999"
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = "'60.64"
$ws.Range('E24').Value = '  -1.98%  '
$ws.Range('D25').Value = "'0.419"
$ws.Range('E25').Value = '  -3.44%  '
$ws.Range('D26').Value = "'0.997"
$ws.Range('D27').Value = '2.683.57'
$ws.Range('E27').Value = '  -4.02%  '
$ws.Range('E28').Value = '  -5.46%  '
$ws.Range('D29').Value = '0.0₃0814'
$ws.Range('E29').Value = '  -5.52%  '
$ws.Range('E30').Value = '  -5.30%  '
$ws.Range('D31').Value = "'0.998"
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').Value = "'6.12"
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('D33').Value = "'18.83"
$ws.Range('E33').Value = '  -3.92%  '
$ws.Range('D34').Value = "'148.98"
$ws.Range('E34').Value = '  -0.64%  '
$ws.Range('E35').Value = '  -5.67%  '
$ws.Range('D36').Value = "'0.949"
$ws.Range('E36').Value = '  +6.69%  '
$ws.Range('D37').Value = "'3.97"
$ws.Range('E37').Value = '  -4.21%  '
$ws.Range('E38').Value = '  -6.22%  '
$ws.Range('D39').Value = "'0.851"
$ws.Range('E39').Value = '  -7.34%  '
$ws.Range('D40').Value = "'36.04"
$ws.Range('E40').Value = '  -2.25%  '
$ws.Range('D41').Value = "'289.50"
$ws.Range('E41').Value = '  -5.63%  '
$ws.Range('E42').Value = '  -5.40%  '
$ws.Range('D43').Value = "'1.39"
$ws.Range('E43').Value = '  -6.99%  '
$ws.Range('D44').Value = "'0.0989"
$ws.Range('E44').Value = '  -3.00%  '
$ws.Range('D45').Value = "'0.996"
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('D46').Value = "'0.609"
$ws.Range('E46').Value = '  -6.46%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'19.24"
$ws.Range('E47').Value = '  -6.57%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').Value = "'0.0536"
$ws.Range('E48').Value = '  -5.51%  '
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('D50').Value = "'0.0228"
$ws.Range('E50').Value = '  -5.50%  '
$ws.Range('E51').Value = '  -8.43%  '
